$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("hindcast")
$ws.Range("B2").Value = 291328965.2223744
$ws.Range("C2").Value = 148779300.5490298
$ws.Range("B3").Value = 118627655.583184
$ws.Range("C3").Value = 60042695.4728971
$ws.Range("B4").Value = 97466474.35908471
$ws.Range("C4").Value = 45843415.20540842
$ws.Range("B5").Value = 113469186.6510352
$ws.Range("C5").Value = 53240208.27381698
$ws.Range("B6").Value = 141101576.4406697
$ws.Range("C6").Value = 64141703.75306258
$ws.Range("B7").Value = 136791840.6839094
$ws.Range("C7").Value = 58381319.42861847
$ws.Range("B8").Value = 113616903.5703495
$ws.Range("C8").Value = 44886926.19054909
$ws.Range("B9").Value = 137544088.0431332
$ws.Range("C9").Value = 61306501.96701577
$ws.Range("B10").Value = 177562095.9139242
$ws.Range("C10").Value = 88494087.23657796
$ws.Range("B11").Value = 197072001.537502
$ws.Range("C11").Value = 102023241.1526247
$ws.Range("B12").Value = 188953745.9961062
$ws.Range("C12").Value = 101054180.1760379
$ws.Range("B13").Value = 175539152.2052628
$ws.Range("C13").Value = 94463219.06450224
$ws.Range("B14").Value = 165576041.1016922
$ws.Range("C14").Value = 89517557.83111714

$ws2 = $wb.Worksheets.Item("condensed")
$ws2.Range("B2").Value = 292770456.2796552
$ws2.Range("C2").Value = 148646848.6585251
$ws2.Range("B3").Value = 119711220.5090676
$ws2.Range("C3").Value = 59850968.41606086
$ws2.Range("B4").Value = 96996904.66162857
$ws2.Range("C4").Value = 45999415.14741379
$ws2.Range("B5").Value = 113754361.1988907
$ws2.Range("C5").Value = 53601793.66675562
$ws2.Range("B6").Value = 142453110.6810186
$ws2.Range("C6").Value = 64408188.21939725
$ws2.Range("B7").Value = 137444135.8373561
$ws2.Range("C7").Value = 58394502.71098745
$ws2.Range("B8").Value = 114653986.0097603
$ws2.Range("C8").Value = 44594793.67013987
$ws2.Range("B9").Value = 138634799.8081673
$ws2.Range("C9").Value = 61526453.42408491
$ws2.Range("B10").Value = 177739222.8826441
$ws2.Range("C10").Value = 88533447.5803074
$ws2.Range("B11").Value = 197737361.8143129
$ws2.Range("C11").Value = 102201838.5866841
$ws2.Range("B12").Value = 189121636.6051421
$ws2.Range("C12").Value = 100805922.1500152
$ws2.Range("B13").Value = 177031362.9710261
$ws2.Range("C13").Value = 94705631.30082849
$ws2.Range("B14").Value = 167766502.971985
$ws2.Range("C14").Value = 89651897.57167193
